$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.918.61'
$ws.Range("D3").Value = '3.684.90'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '2.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +11.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '235.17'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '653.39'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.433'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.35%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = '3.684.00'
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000308'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +14.87%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.27'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.206'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.78'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '4.374.56'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = '96.617.07'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.76'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.54%  '
$ws.Range("D19").Value = '3.694.82'
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.89'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.50'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.533'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '514.45'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.41'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.96%  '
$ws.Range("E25").Value = '  +6.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.85'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '109.41'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.201'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +19.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.30'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.48'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.98'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.186'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.80'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.41'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.587'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '628.02'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.66'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.46%  '
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.490'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +7.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.73'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.00'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.15'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.949'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.03%  '
$ws.Range("E47").Value = '  -3.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.36'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.60'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.62'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.31'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.46%  '
